$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the empty "Sheet2" tab. Excel keeps "Sheet3" (name/sheetId kept),
#    and on save it will be re-keyed onto the freed relationship id / the
#    now-available "sheet2.xml" part name.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

$s1 = $wb.Worksheets.Item("Sheet1")
$s3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 2. Fix the RD-amount label typo on Sheet1 ("RD AMONUT" -> "RD AMOUNT").
#    Doing this first means the new corrected string becomes the first fresh
#    shared-string entry, matching the authoritative index ordering.
# ---------------------------------------------------------------------------
$s1.Range("B25").Value = "RD AMOUNT"

# ---------------------------------------------------------------------------
# 3. Re-point Sheet1!E5 at the (erroneous) Sheet3!EC5 cell reference and let
#    the whole dependent chain (C6:F21, shared formulas) recompute.
# ---------------------------------------------------------------------------
$s1.Range("E5").Formula = "=SUM((Sheet3!EC5+D5)*7.1%/12)"

# ---------------------------------------------------------------------------
# 4. Populate "Sheet3" with the recurring-deposit table (a copy of Sheet1's
#    first table, shifted up one row, columns B:F).
# ---------------------------------------------------------------------------
$s3.Range("C2").Value = "OPENING BALNACE"
$s3.Range("D2").Value = "PRINICPAL AMOUNT"
$s3.Range("E2").Value = "INTEREST REATE 7.10%"
$s3.Range("F2").Value = "PRINCIPAL AMOUNT WITH RATE "
$s3.Range("B2").Value = "MONTH"

# Month numbers 1..12 in B3:B14 (plain literal values, not formulas).
for ($i = 0; $i -lt 12; $i++) {
    $s3.Cells.Item(3 + $i, 2).Value = $i + 1
}

# Principal deposited every month (D3:D14).
$s3.Range("D3:D14").Value = 1000

# Opening balance chain (C4:C14) + interest (E3:E14) + running total (F3:F14).
$s3.Range("E3").Formula = "=SUM(D3*7.1%/12)"
$s3.Range("F3").Formula = "=SUM(D3+E3)"

$s3.Range("C4").Formula = "=F3"
$s3.Range("E4").Formula = "=SUM((C4+D4)*7.1%/12)"
$s3.Range("F4").Formula = "=SUM(C4+D4+E4)"

$s3.Range("C5:C14").Formula = "=F4"
$s3.Range("E5:E14").Formula = "=SUM((C5+D5)*7.1%/12)"
$s3.Range("F5:F14").Formula = "=SUM(C5+D5+E5)"

# Number formatting (2 decimal places) on the computed columns only.
$s3.Range("E3:E14").NumberFormat = "0.00"
$s3.Range("F3:F14").NumberFormat = "0.00"
$s3.Range("C4:C14").NumberFormat = "0.00"

# Trailing blank-label cell.
$s3.Range("E19").Value = " "

# Column widths (best-fit approximation) for the populated columns.
$s3.Range("B2:F14").EntireColumn.AutoFit() | Out-Null

# Print setup to match Sheet1's page setup.
$s3.PageSetup.PaperSize = 9
$s3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Selections / active tab: Sheet1 no longer the displayed tab; Sheet3
#    becomes the active (second) tab with its own selection.
# ---------------------------------------------------------------------------
$s1.Range("E5").Select() | Out-Null
$s3.Activate() | Out-Null
$s3.Range("J10").Select() | Out-Null

Write-Output "done"
